$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value2 = "20.196.55"
$ws.Range("E2").Value2 = "  +2.09%  "
$ws.Range("D3").Value2 = "1.434.94"
$ws.Range("E3").Value2 = "  +3.51%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value2 = "1.008"
$ws.Range("E4").Value2 = "  +0.79%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value2 = "0.9128"
$ws.Range("E5").Value2 = "  -8.88%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value2 = "276.62"
$ws.Range("E6").Value2 = "  +3.21%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value2 = "0.3637"
$ws.Range("E7").Value2 = "  +0.31%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value2 = "0.3093"
$ws.Range("E8").Value2 = "  +1.77%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value2 = "39.02"
$ws.Range("E9").Value2 = "  -0.25%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value2 = "1.016"
$ws.Range("E10").Value2 = "  +3.98%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value2 = "0.06508"
$ws.Range("E11").Value2 = "  +1.38%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value2 = "1.003"
$ws.Range("E12").Value2 = "  +0.29%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value2 = "5.337"
$ws.Range("E13").Value2 = "  +0.80%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value2 = "17.48"
$ws.Range("E14").Value2 = "  +4.28%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value2 = "6.041"
$ws.Range("E15").Value2 = "  -0.38%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value2 = "0.00001011"
$ws.Range("E16").Value2 = "  +1.71%  "
$ws.Range("D17").Value2 = "1.438.25"
$ws.Range("E17").Value2 = "  +3.61%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value2 = "0.9410"
$ws.Range("E18").Value2 = "  -6.05%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value2 = "0.05638"
$ws.Range("E19").Value2 = "  +0.37%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value2 = "67.72"
$ws.Range("E20").Value2 = "  -2.99%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value2 = "5.347"
$ws.Range("E21").Value2 = "  -2.87%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value2 = "14.26"
$ws.Range("E22").Value2 = "  -1.13%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value2 = "10.76"
$ws.Range("E23").Value2 = "  +1.98%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value2 = "2.248"
$ws.Range("E24").Value2 = "  +0.44%  "
$ws.Range("D25").Value2 = "20.269.89"
$ws.Range("E25").Value2 = "  +2.46%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value2 = "2.133"
$ws.Range("E26").Value2 = "  -1.97%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value2 = "136.52"
$ws.Range("E27").Value2 = "  -0.23%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value2 = "16.86"
$ws.Range("E28").Value2 = "  +1.99%  "
$ws.Range("D29").Value2 = "1.589.54"
$ws.Range("E29").Value2 = "  +3.07%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value2 = "109.49"
$ws.Range("E30").Value2 = "  +1.74%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value2 = "3.900"
$ws.Range("E31").Value2 = "  +1.68%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value2 = "0.7914"
$ws.Range("E32").Value2 = "  -1.16%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value2 = "4.803"
$ws.Range("E33").Value2 = "  -7.96%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value2 = "0.07660"
$ws.Range("E34").Value2 = "  +1.11%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value2 = "0.05886"
$ws.Range("E35").Value2 = "  +4.64%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value2 = "1.446"
$ws.Range("E36").Value2 = "  +11.75%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value2 = "1.136"
$ws.Range("E37").Value2 = "  +8.42%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value2 = "4.600"
$ws.Range("E38").Value2 = "  -1.69%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value2 = "0.01979"
$ws.Range("E39").Value2 = "  -2.34%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value2 = "10.13"
$ws.Range("E40").Value2 = "  +1.20%  "
$ws.Range("B41").Value2 = "Algorand"
$ws.Range("C41").Value2 = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value2 = "0.1833"
$ws.Range("E41").Value2 = "  -2.00%  "
$ws.Range("B42").Value2 = "Frax"
$ws.Range("C42").Value2 = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value2 = "0.9216"
$ws.Range("E42").Value2 = "  -7.91%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value2 = "6.996"
$ws.Range("E43").Value2 = "  -14.65%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value2 = "3.500"
$ws.Range("E44").Value2 = "  +1.40%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value2 = "0.5198"
$ws.Range("E45").Value2 = "  +0.20%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value2 = "11.91"
$ws.Range("E46").Value2 = "  -0.92%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value2 = "117.52"
$ws.Range("E47").Value2 = "  +7.69%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value2 = "0.5089"
$ws.Range("E48").Value2 = "  +2.16%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value2 = "1.747"
$ws.Range("E49").Value2 = "  +1.19%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value2 = "0.06313"
$ws.Range("E50").Value2 = "  +4.34%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value2 = "0.9868"
$ws.Range("E51").Value2 = "  -1.32%  "
